$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 4: a paper from TUD, assigned to 张睿 (Rui) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "FastLan: Improving Performance of Software Transactional Memory for Low Thread Counts"
$ws.Range("C4").Value = "Jons-Tobias Wamhoff, etc."
$ws.Range("D4").Value = "PPoPP 2013"

# Force the date-looking text to stay literal text instead of being parsed as a date
# (mirrors typing a leading apostrophe in Excel to force text entry).
$ws.Range("E4").Value = "'2013.2.23"

$ws.Range("F4").Value = "张睿"
$ws.Range("G4").Value = "张睿"

# --- New column H: remarks ---
$ws.Range("H1").Value = "备注"
$ws.Range("H4").Value = "TUD的一个做并行的组的成果，关注一下。"

# Widen column H to fit the new remark text.
$ws.Columns.Item(8).ColumnWidth = 39

# Match the post-edit selection.
$ws.Range("H5").Select()
